$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are stored as literal text (coinranking.com scrape output),
# so force each target cell to Text format before writing the new reading. This
# keeps e.g. "320.24" and "7.49%" as text instead of being auto-converted to a
# number/percentage by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "320.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "7.49%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "48.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "16.74%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.279"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.47%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08095"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.603"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.22%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.656"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.07%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.201"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "30.81%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1318"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "11.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1943"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09436"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.39%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04507"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "9.82%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.06%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001323"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.42%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005935"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.25%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.61%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.437"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.48%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3392"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.91%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.236"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.32%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.99%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3059"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04291"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.90%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001310"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.40%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "8.24%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001352"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.96%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003545"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.81%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02677"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "11.34%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05569"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006311"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.09%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007681"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.98%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1437"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007707"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.12%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.27%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.86%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006995"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.15%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "34.76%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004006"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.68%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
